{"js": "// Add a new 4th row to the table, mirroring the existing rows' layout.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst newRow = [\n  \"4)\",\n  \"szer. pchor.\",\n  \"dsa\",\n  \"DSA\",\n  \"w dn. 26 - 26.10.2020 r.\",\n  \"do m.\",\n  \"dsa\",\n];\n\ntable.addRows(\"End\", 1, [newRow]);\nawait context.sync();\n", "ps1": "# Add a new 4th row to the end of the (only) table, mirroring the\n# existing rows' layout/columns.\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newRow = $t.Rows.Add()\n$i = $newRow.Index\n\n$t.Cell($i, 1).Range.Text = \"4)\"\n$t.Cell($i, 2).Range.Text = \"szer. pchor.\"\n$t.Cell($i, 3).Range.Text = \"dsa\"\n$t.Cell($i, 4).Range.Text = \"DSA\"\n$t.Cell($i, 5).Range.Text = \"w dn. 26 - 26.10.2020 r.\"\n$t.Cell($i, 6).Range.Text = \"do m.\"\n$t.Cell($i, 7).Range.Text = \"dsa\"\n"}
